# Commit: "@tongnd : commit for Vy"
#
# Add a new task line to the "IT" sheet noting that a source/environment
# needs to be created for sukien.talenttech6.vn, and make the "IT" sheet
# the active/selected sheet (it had previously been left on "Test").

$wb = $excel.ActiveWorkbook

$itSheet = $wb.Worksheets.Item("IT")

# New row 6 on the IT sheet with the note about the new task.
$itSheet.Range("A6").Value = "Tạo source cho sukien.talenttech6.vn"

# Make "IT" the active sheet (moves tabSelected from "Test" to "IT")
# and leave the same kind of "last looked at" cell selected as in the
# saved file.
$itSheet.Activate()
$itSheet.Range("G31").Select()
